$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the data range so numeric-looking strings
# (e.g. "604.43", "0.0000256") are not auto-converted to numbers,
# matching the source data which stores these as literal text.
$ws.Range("B2:E51").NumberFormat = "@"

# Apply updated crypto price/volume data as per commit
$ws.Range("D2").Value = '69.553.01'
$ws.Range("E2").Value = '  +1.05%  '
$ws.Range("D3").Value = '3.881.18'
$ws.Range("E3").Value = '  +0.32%  '
$ws.Range("E4").Value = '  -0.16%  '
$ws.Range("D5").Value = '604.43'
$ws.Range("E5").Value = '  +0.67%  '
$ws.Range("D6").Value = '169.73'
$ws.Range("E6").Value = '  +4.52%  '
$ws.Range("D7").Value = '3.876.91'
$ws.Range("E7").Value = '  +0.21%  '
$ws.Range("E8").Value = '  +0.12%  '
$ws.Range("D9").Value = '0.535'
$ws.Range("E9").Value = '  +0.98%  '
$ws.Range("D10").Value = '0.168'
$ws.Range("E10").Value = '  -0.15%  '
$ws.Range("E11").Value = '  +0.90%  '
$ws.Range("D12").Value = '0.468'
$ws.Range("E12").Value = '  +2.06%  '
$ws.Range("D13").Value = '0.0000256'
$ws.Range("E13").Value = '  +5.08%  '
$ws.Range("D14").Value = '38.29'
$ws.Range("E14").Value = '  +3.75%  '
$ws.Range("D15").Value = '4.535.92'
$ws.Range("E15").Value = '  +0.37%  '
$ws.Range("D16").Value = '3.886.25'
$ws.Range("E16").Value = '  +0.40%  '
$ws.Range("D17").Value = '69.532.94'
$ws.Range("E17").Value = '  +0.78%  '
$ws.Range("D18").Value = '18.73'
$ws.Range("E18").Value = '  +9.54%  '
$ws.Range("D19").Value = '7.66'
$ws.Range("E19").Value = '  +1.49%  '
$ws.Range("E20").Value = '  -0.84%  '
$ws.Range("D21").Value = '11.18'
$ws.Range("E21").Value = '  -1.56%  '
$ws.Range("D22").Value = '490.53'
$ws.Range("E22").Value = '  +1.39%  '
$ws.Range("D23").Value = '0.748'
$ws.Range("E23").Value = '  +4.11%  '
$ws.Range("D24").Value = '0.0000167'
$ws.Range("E24").Value = '  +3.62%  '
$ws.Range("D25").Value = '85.28'
$ws.Range("E25").Value = '  +1.57%  '
$ws.Range("D26").Value = '2.31'
$ws.Range("E26").Value = '  +3.28%  '
$ws.Range("D27").Value = '12.31'
$ws.Range("E27").Value = '  +1.91%  '
$ws.Range("D28").Value = '10.17'
$ws.Range("E28").Value = '  +2.48%  '
$ws.Range("E29").Value = '  +0.11%  '
$ws.Range("D30").Value = '2.98'
$ws.Range("E30").Value = '  +0.92%  '
$ws.Range("D31").Value = '2.43'
$ws.Range("E31").Value = '  +2.55%  '
$ws.Range("D32").Value = '4.032.17'
$ws.Range("E32").Value = '  +0.29%  '
$ws.Range("E33").Value = '  -1.29%  '
$ws.Range("D34").Value = '31.92'
$ws.Range("E34").Value = '  -1.13%  '
$ws.Range("D35").Value = '3.845.11'
$ws.Range("E35").Value = '  +0.75%  '
$ws.Range("E36").Value = '  +0.39%  '
$ws.Range("D37").Value = '6.12'
$ws.Range("E37").Value = '  +4.23%  '
$ws.Range("D38").Value = '1.04'
$ws.Range("E38").Value = '  +0.41%  '
$ws.Range("D39").Value = '0.142'
$ws.Range("E39").Value = '  +1.12%  '
$ws.Range("D40").Value = '3.30'
$ws.Range("E40").Value = '  +11.31%  '
$ws.Range("D41").Value = '0.999'
$ws.Range("E41").Value = '  -0.11%  '
$ws.Range("D42").Value = '0.328'
$ws.Range("E42").Value = '  +3.00%  '
$ws.Range("D43").Value = '2.12'
$ws.Range("E43").Value = '  +6.78%  '
$ws.Range("D44").Value = '438.27'
$ws.Range("E44").Value = '  +1.06%  '
$ws.Range("D45").Value = '48.16'
$ws.Range("E45").Value = '  -0.68%  '
$ws.Range("D46").Value = '8.72'
$ws.Range("E46").Value = '  +3.94%  '
$ws.Range("D48").Value = '0.0369'
$ws.Range("E48").Value = '  +3.03%  '
$ws.Range("D49").Value = '143.59'
$ws.Range("E49").Value = '  +0.12%  '
$ws.Range("B50").Value = 'FLOKI'
$ws.Range("C50").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range("D50").Value = '0.000271'
$ws.Range("E50").Value = '  +19.34%  '
$ws.Range("B51").Value = 'Arweave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D51").Value = '40.13'
$ws.Range("E51").Value = '  +3.96%  '

# Restore default (Normal) style on the data range so no stray
# number-format/style attribute is left behind on the cells.
$ws.Range("B2:E51").Style = "Normal"
